$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bulk numeric corrections across rows 128-219 (N/O/C/D columns) ---
$cellValues = @{
  "O128" = 130
  "O129" = 129
  "O130" = 125
  "O131" = 119
  "O132" = 110
  "O133" = 93
  "O134" = 84
  "O135" = 75
  "O136" = 71
  "O137" = 53
  "O138" = 75
  "O139" = 80
  "O140" = 83
  "O141" = 85
  "O142" = 82
  "O143" = 83
  "O144" = 82
  "O145" = 80
  "O146" = 63
  "O147" = 66
  "O148" = 48
  "O149" = 50
  "O150" = 57
  "O151" = 56
  "O162" = 157
  "O163" = 170
  "O164" = 185
  "O165" = 179
  "O166" = 167
  "O167" = 194
  "O168" = 193
  "O169" = 209
  "O170" = 216
  "O171" = 231
  "O172" = 221
  "O173" = 196
  "O174" = 242
  "O175" = 259
  "O176" = 275
  "O177" = 285
  "O178" = 279
  "O179" = 286
  "O180" = 361
  "O181" = 366
  "O182" = 373
  "O183" = 414
  "O184" = 434
  "O185" = 460
  "O186" = 446
  "O187" = 427
  "O188" = 358
  "O189" = 317
  "O190" = 297
  "O191" = 261
  "O192" = 223
  "O193" = 233
  "O194" = 242
  "O195" = 213
  "O196" = 214
  "O197" = 225
  "O198" = 213
  "O199" = 227
  "O200" = 234
  "O201" = 258
  "O202" = 296
  "O203" = 290
  "O204" = 332
  "O205" = 386
  "O206" = 409
  "O207" = 428
  "O208" = 421
  "O209" = 386
  "O210" = 400
  "O211" = 405
  "C212" = 10
  "N212" = 109
  "O212" = 429
  "N213" = 107
  "O213" = 448
  "N214" = 95
  "O214" = 394
  "N215" = 94
  "O215" = 351
  "O216" = 312
  "C217" = 18
  "O217" = 302
  "C218" = 21
  "N218" = 105
  "O218" = 279
  "C219" = 15
  "D219" = 1
  "N219" = 116
  "O219" = 349
}

foreach ($addr in $cellValues.Keys) {
  $ws.Range($addr).Value2 = $cellValues[$addr]
}

# --- New row 220 (previously blank placeholder row) ---
$ws.Range("C220").Value2 = 0
$ws.Range("D220").Value2 = 0
$ws.Range("E220").Value2 = 1
$ws.Range("F220").Value2 = 1
$ws.Range("G220").Value2 = 8
$ws.Range("I220").Value2 = 0
$ws.Range("L220").Value2 = "0"
$ws.Range("M220").Value2 = "0"
$ws.Range("N220").Value2 = 98
$ws.Range("O220").Value2 = 330
$ws.Range("P220").Value2 = "na"

# --- Selection moves to N222 (matches the saved view state) ---
$ws.Range("N222").Select()
